$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

$ws1.Range("B1").Value = "4220-RBI-EI-DB-DL-REC-RNI-FEE+INT-FFConMONonLASTSUN-FIFC-1-FFROP-DLY-FIFR-1-MD-TR1-OT1st"
$ws2.Range("B1").Value = "4220-RBI-EI-DB-DL-REC-RNI-FEE+INT-FFConMONonLASTSUN-FIFC-1-FFROP-DLY-FIFR-1-MD-TR1-OT1st"
$ws1.Range("B2").Value = "422t"
$ws1.Range("B13").Select()
$ws2.Activate()
